{"js": "// G ajoute des truc\n// Adds three new paragraphs (two blank + one note to the team) right\n// before the very last (empty) paragraph of the document, i.e. just\n// after the \"Plan d'action\" bullet item. The \"_GoBack\" bookmark that\n// used to sit on the bullet item is moved along to the new last\n// paragraph, matching how Word re-anchors it to the most recent edit\n// location.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document always ends with a trailing empty paragraph; insert the\n// new content right before it so the new paragraphs inherit its plain\n// (non-list) formatting instead of the bulleted-list formatting of the\n// \"Plan d'action\" item above it.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Drop the existing \"_GoBack\" bookmark - it will be re-created at the\n// new end of the document below.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst noteParagraph = lastParagraph.insertParagraph(\n  \"WSH oublie pas de modifier le bas !!\",\n  Word.InsertLocation.before\n);\nnoteParagraph.insertParagraph(\"\", Word.InsertLocation.before);\nnoteParagraph.insertParagraph(\"\", Word.InsertLocation.before);\n\n// Re-create \"_GoBack\" on the (still) final paragraph of the document.\nlastParagraph.getRange().insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# G ajoute des truc\n# Adds three new paragraphs (two blank + one note to the team) right\n# before the very last (empty) paragraph of the document, i.e. just\n# after the \"Plan d'action\" bullet item. The \"_GoBack\" bookmark that\n# used to sit on the bullet item is moved along to the new last\n# paragraph, matching how Word re-anchors it to the most recent edit\n# location.\n\n$d = $word.ActiveDocument\n\n# Drop the existing \"_GoBack\" bookmark - it will be re-created at the\n# new end of the document below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# The document always ends with a trailing empty paragraph; insert the\n# new content right before it so the new paragraphs inherit its plain\n# (non-list) formatting instead of the bulleted-list formatting of the\n# \"Plan d'action\" item above it.\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$rng = $lastParagraph.Range\n$rng.InsertParagraphBefore()\n$rng.InsertParagraphBefore()\n$rng.InsertParagraphBefore()\n\n$notePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n$notePara.Range.Text = \"WSH oublie pas de modifier le bas !!\"\n\n# Re-create \"_GoBack\" on the (still) final paragraph of the document.\n$finalParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$d.Bookmarks.Add(\"_GoBack\", $finalParagraph.Range)\n"}
